# Adding Empty key validation
# This script reorganizes the exception table: the "InvalidKeyException /
# Parameters missing / When IV is not provided in CBC" row moves from row 12
# down to row 19 (after the AEADBadTagException block), the BadPaddingException
# and AEADBadTagException blocks shift up accordingly, and a brand new row 20
# is appended containing the new SAF-015 "Unsecure Key" validation entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final desired contents for rows 12..20, columns A..E ($null = blank cell)
$rows = @(
    @("BadPaddingException", "Given final block not properly padded. Such issues can arise if a bad key is used during decryption", "Wrong IV", "SAF-010", "Either the Mode/Key/IV used for encryption was different than provided for decryption"),
    @("BadPaddingException", "Given final block not properly padded. Such issues can arise if a bad key is used during decryption", "Wrong Key", "SAF-010", "Either the Mode/Key/IV used for encryption was different than provided for decryption"),
    @("BadPaddingException", "Given final block not properly padded. Such issues can arise if a bad key is used during decryption", "Wrong Mode (Encrypted in ECB, Decrypting in CBC or vice versa)", "SAF-010", "Either the Mode/Key/IV used for encryption was different than provided for decryption"),
    @("BadPaddingException", "Given final block not properly padded. Such issues can arise if a bad key is used during decryption", "Wrong Padding (Different for encryption, and different for decryption)", "SAF-010", "Either the Mode/Key/IV used for encryption was different than provided for decryption"),
    @("AEADBadTagException", "Tag mismatch", "Wrong IV provided when using GCM Mode", "SAF-002", "Either the Key/IV/Associated Data used for encryption was different than provided for decryption"),
    @("AEADBadTagException", "Tag mismatch", "Wrong Key provided when using GCM Mode", "SAF-002", "Either the Key/IV/Associated Data used for encryption was different than provided for decryption"),
    @("AEADBadTagException", "Tag mismatch", "Wrong IV provided when using GCM Mode", "SAF-002", "Either the Key/IV/Associated Data used for encryption was different than provided for decryption"),
    @("InvalidKeyException", "Parameters missing", "When IV is not provided in CBC", $null, "Compile Time: Doesn’t Allow to use without IV"),
    @($null, $null, $null, "SAF-015", "Compile Time: Unsecure Key! Key is initialized with all zeros")
)

$startRow = 12
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $vals = $rows[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $col = $c + 1
        $v = $vals[$c]
        if ($null -eq $v) {
            $ws.Cells.Item($r, $col).Value = $null
        } else {
            $ws.Cells.Item($r, $col).Value = $v
        }
    }
}

# Update the selection / view to match the committed state
$ws.Range("C20").Select()
